$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-numeric-looking text columns (B, C, E) can be set directly.
$ws.Range('E2').Value = '  +8.44%  '
$ws.Range('E3').Value = '  +7.00%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  +4.36%  '
$ws.Range('E8').Value = '  +8.82%  '
$ws.Range('E9').Value = '  +5.37%  '
$ws.Range('E10').Value = '  +5.34%  '
$ws.Range('E11').Value = '  +8.82%  '
$ws.Range('E12').Value = '  +7.95%  '
$ws.Range('E13').Value = '  +6.05%  '
$ws.Range('E14').Value = '  +6.54%  '
$ws.Range('E15').Value = '  +5.07%  '
$ws.Range('E16').Value = '  +3.53%  '
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('E18').Value = '  +4.52%  '
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('E20').Value = '  +5.86%  '
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').Value = '  +8.37%  '
$ws.Range('E23').Value = '  +6.09%  '
$ws.Range('E24').Value = '  +4.19%  '
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('E26').Value = '  +6.02%  '
$ws.Range('E27').Value = '  +2.48%  '
$ws.Range('E28').Value = '  +4.68%  '
$ws.Range('E29').Value = '  +7.69%  '
$ws.Range('E30').Value = '  +8.23%  '
$ws.Range('E31').Value = '  +4.37%  '
$ws.Range('E32').Value = '  +10.12%  '
$ws.Range('E33').Value = '  +3.57%  '
$ws.Range('E34').Value = '  +12.58%  '
$ws.Range('E35').Value = '  +3.27%  '
$ws.Range('E36').Value = '  +5.84%  '
$ws.Range('E38').Value = '  +6.60%  '
$ws.Range('E39').Value = '  +6.38%  '
$ws.Range('E40').Value = '  +5.06%  '
$ws.Range('E42').Value = '  +8.46%  '
$ws.Range('E43').Value = '  +4.84%  '
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('E46').Value = '  +32.67%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E47').Value = '  +7.54%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('E48').Value = '  +6.34%  '
$ws.Range('E49').Value = '  +7.04%  '
$ws.Range('E50').Value = '  +12.58%  '

# Column D holds numeric-looking strings (prices). Force text format
# so Excel doesn't coerce them to floats (which would drop trailing zeros
# like '1.040' -> 1.04, or reformat '29.739.84' style multi-dot prices).
$dRange = $ws.Range('D2:D51')
$dRange.NumberFormat = '@'
$ws.Range('D2').Value = '29.739.84'
$ws.Range('D3').Value = '1.948.67'
$ws.Range('D5').Value = '342.09'
$ws.Range('D7').Value = '0.4780'
$ws.Range('D8').Value = '0.4145'
$ws.Range('D9').Value = '48.29'
$ws.Range('D10').Value = '0.08257'
$ws.Range('D11').Value = '1.040'
$ws.Range('D12').Value = '22.66'
$ws.Range('D13').Value = '1.941.45'
$ws.Range('D14').Value = '6.210'
$ws.Range('D15').Value = '7.422'
$ws.Range('D16').Value = '92.30'
$ws.Range('D18').Value = '0.00001064'
$ws.Range('D19').Value = '0.06682'
$ws.Range('D20').Value = '18.05'
$ws.Range('D22').Value = '29.701.54'
$ws.Range('D23').Value = '5.606'
$ws.Range('D24').Value = '11.26'
$ws.Range('D25').Value = '2.286'
$ws.Range('D26').Value = '2.173.42'
$ws.Range('D27').Value = '161.12'
$ws.Range('D28').Value = '20.20'
$ws.Range('D30').Value = '5.689'
$ws.Range('D31').Value = '122.66'
$ws.Range('D33').Value = '0.09645'
$ws.Range('D34').Value = '1.479'
$ws.Range('D35').Value = '3.683'
$ws.Range('D36').Value = '5.514'
$ws.Range('D37').Value = '0.06321'
$ws.Range('D38').Value = '0.02330'
$ws.Range('D39').Value = '8.595'
$ws.Range('D40').Value = '1.197'
$ws.Range('D41').Value = '0.6121'
$ws.Range('D42').Value = '10.75'
$ws.Range('D43').Value = '0.1901'
$ws.Range('D45').Value = '1.273'
$ws.Range('D46').Value = '2.390'
$ws.Range('D47').Value = '12.61'
$ws.Range('D48').Value = '0.5727'
$ws.Range('D49').Value = '2.003'
$ws.Range('D50').Value = '0.07397'
$ws.Range('D51').Value = '114.20'
$dRange.Style = 'Normal'

